# FanPowerAllowances-T24N_2022.xlsx edit
#
# The "FanPwrIdxAdj" table (rows 1-11) previously had a single "ERV"
# (Energy Recovery) column.  It is split into four columns:
#   Exhaust Systems Base Allowance / Supply-ERV / Return-ERV / Return-Filter
# which pushes the trailing "SZVAV" note column from K to N.
#
# The lower "BaseFanPwrIdx" table (rows 13-19) lives in the same
# columns (A-L) further down the sheet and must NOT move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert three whole columns at J:L.  (The runtime's Range.Insert
#    only supports whole-column shifts reliably, so we do a full
#    column insert and then slide the lower table back into place.)
# ---------------------------------------------------------------
$ws.Columns("J:L").Insert()

# The column insert also pushed the second ("BaseFanPwrIdx") table's
# J:L columns over to M:O -- move that block back to J:L so only the
# first table's columns actually grew.
$ws.Range("M14:O18").Cut($ws.Range("J14:L18"))

# ---------------------------------------------------------------
# 2. Clear the inherited (old-J / ERV) number formatting from the new
#    J:M cells of the top table before re-populating them.
# ---------------------------------------------------------------
$ws.Range("J3:M10").ClearFormats()

# ---------------------------------------------------------------
# 3. Row 3 header text - reuse the wrapped/vertical-top style that the
#    rest of row 3 (D3:I3) already has.
# ---------------------------------------------------------------
$ws.Range("I3").Copy()
$ws.Range("J3:M3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("J3").Value = "Exhaust Systems Base Allowance"
$ws.Range("K3").Value = "Supply - Energy Recovery (Enthalpy Recovery Ratio ? 0.60 and <0.65)  "
$ws.Range("L3").Value = "Return - Energy Recovery (Enthalpy Recovery Ratio ? 0.60 and <0.65)  "
$ws.Range("M3").Value = "Return - Filter (any MERV value)"

# ---------------------------------------------------------------
# 4. Row 4 (variable-name row).  New header cells get the wrap-text
#    style used elsewhere in the sheet (e.g. A14).
# ---------------------------------------------------------------
$ws.Range("J4").Value = "ExhBaseAllow"
$ws.Range("K4").Value = "ERVSupply"
$ws.Range("L4").Value = "ERVReturn"
$ws.Range("M4").Value = "RetFilter"
$ws.Range("J4:M4").WrapText = $true

$ws.Rows("4").RowHeight = 30

# ---------------------------------------------------------------
# 5. Numeric data rows 5-10 (new J:M values).  These stay in the
#    default/General number format (no 0.000 style), matching the
#    target workbook.
# ---------------------------------------------------------------
$ws.Range("J5").Value = 0.221
$ws.Range("K5").Value = 0.184
$ws.Range("L5").Value = 0.19
$ws.Range("M5").Value = 0.046

$ws.Range("J6").Value = 0.246
$ws.Range("K6").Value = 0.155
$ws.Range("L6").Value = 0.163
$ws.Range("M6").Value = 0.041

$ws.Range("J7").Value = 0.236
$ws.Range("K7").Value = 0.144
$ws.Range("L7").Value = 0.146
$ws.Range("M7").Value = 0.036

$ws.Range("J8").Value = 0.186
$ws.Range("K8").Value = 0.19
$ws.Range("L8").Value = 0.191
$ws.Range("M8").Value = 0.046

$ws.Range("J9").Value = 0.184
$ws.Range("K9").Value = 0.163
$ws.Range("L9").Value = 0.166
$ws.Range("M9").Value = 0.041

$ws.Range("J10").Value = 0.19
$ws.Range("K10").Value = 0.146
$ws.Range("L10").Value = 0.148
$ws.Range("M10").Value = 0.036

# ---------------------------------------------------------------
# 6. Column widths: D:I stay width 20; new J:N columns get the wider
#    width used for the new data block.
# ---------------------------------------------------------------
$ws.Columns("J:N").ColumnWidth = 21.86

# ---------------------------------------------------------------
# 7. Misc header/row-height/selection housekeeping to mirror the
#    authored workbook.
# ---------------------------------------------------------------
$ws.Rows("3").RowHeight = 90
$ws.Range("M15").Select()
